$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.710.08"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "3.088.68"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'575.92"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").Value = "'177.21"
$ws.Range("E6").Value = "  +2.61%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.086.42"
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("D9").Value = "'0.513"
$ws.Range("E9").Value = "  -1.34%  "
$ws.Range("D10").Value = "'6.35"
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("D11").Value = "'0.151"
$ws.Range("E11").Value = "  -1.86%  "
$ws.Range("D12").Value = "'0.465"
$ws.Range("E12").Value = "  -2.84%  "
$ws.Range("D13").Value = "'0.0000239"
$ws.Range("E13").Value = "  -3.24%  "
$ws.Range("D14").Value = "'35.88"
$ws.Range("E14").Value = "  -2.45%  "
$ws.Range("D15").Value = "'0.122"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").Value = "3.596.78"
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("D17").Value = "66.626.78"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").Value = "'6.96"
$ws.Range("E18").Value = "  -1.77%  "
$ws.Range("D19").Value = "'16.79"
$ws.Range("E19").Value = "  +2.16%  "
$ws.Range("D20").Value = "3.083.88"
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("D21").Value = "'480.07"
$ws.Range("E21").Value = "  -2.08%  "
$ws.Range("D22").Value = "'7.73"
$ws.Range("E22").Value = "  -2.41%  "
$ws.Range("D23").Value = "'0.686"
$ws.Range("E23").Value = "  -2.42%  "
$ws.Range("D24").Value = "'83.20"
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("D25").Value = "'12.60"
$ws.Range("E25").Value = "  -4.17%  "
$ws.Range("D26").Value = "'2.21"
$ws.Range("E26").Value = "  -3.07%  "
$ws.Range("D27").Value = "'10.08"
$ws.Range("E27").Value = "  -3.63%  "
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").Value = "'7.94"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("D30").Value = "'2.28"
$ws.Range("E30").Value = "  -3.74%  "
$ws.Range("D31").Value = "'2.59"
$ws.Range("E31").Value = "  -2.68%  "
$ws.Range("D32").Value = "'27.90"
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("E33").Value = "  -2.10%  "
$ws.Range("D34").Value = "0.0₃0936"
$ws.Range("E34").Value = "  -0.68%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "'48.09"
$ws.Range("E36").Value = "  +2.50%  "
$ws.Range("D37").Value = "'5.56"
$ws.Range("E37").Value = "  -4.85%  "
$ws.Range("E38").Value = "  -3.77%  "
$ws.Range("D39").Value = "'48.94"
$ws.Range("E39").Value = "  -2.16%  "
$ws.Range("D40").Value = "'0.308"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("D42").Value = "'1.98"
$ws.Range("E42").Value = "  -2.75%  "
$ws.Range("D43").Value = "'8.30"
$ws.Range("E43").Value = "  -1.81%  "
$ws.Range("D44").Value = "'2.67"
$ws.Range("E44").Value = "  +3.44%  "
$ws.Range("D45").Value = "2.780.68"
$ws.Range("E45").Value = "  -0.79%  "
$ws.Range("D46").Value = "'371.07"
$ws.Range("E46").Value = "  -3.13%  "
$ws.Range("D47").Value = "'135.41"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("E48").Value = "  -2.39%  "
$ws.Range("D50").Value = "'24.66"
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("D51").Value = "'2.21"
$ws.Range("E51").Value = "  +0.79%  "
